$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Column A width (closest achievable snap to target 22.1640625 via this engine's 1/6-grid rounding)
$ws.Columns.Item(1).ColumnWidth = 21.3

# --- New row 47: transposed table header ---
$ws.Range("B35:D35").Copy()
$ws.Range("B47:D47").PasteSpecial(-4122)
$ws.Range("B47").Value = 'Index files'
$ws.Range("C47").Value = 'Raw Shortreads'
$ws.Range("D47").Value = 'processed Shortreads'

# --- Column A labels (48:53), format copied from a plain text cell (A33) ---
$ws.Range("A33").Copy()
$ws.Range("A48:A53").PasteSpecial(-4122)
$ws.Range("A48").Value = 'Size'
$ws.Range("A49").Value = 'Cyder to QB(s)'
$ws.Range("A50").Value = 'Cyder to Ranger(S)'
$ws.Range("A51").Value = 'Cyder to QB & QB to Ranger(s)'
$ws.Range("A52").Value = 'QB to Ranger(s)'
$ws.Range("A53").Value = 'Min Total time(s)'

# --- Data cells B:D 48:53, format copied from B35:D37 (style s=7) ---
$ws.Range("B35:D37").Copy()
$ws.Range("B48:D50").PasteSpecial(-4122)
$ws.Range("B51:D53").PasteSpecial(-4122)

$ws.Range("B48").Value = 127
$ws.Range("C48").Value = 9
$ws.Range("D48").Value = 3.81

$ws.Range("B49").Value = 1814.29
$ws.Range("C49").Value = 128.57
$ws.Range("D49").Value = '               N/A'

$ws.Range("B50").Value = 21166.67
$ws.Range("C50").Value = '      N/A'
$ws.Range("D50").Value = '       N/A'

$ws.Range("B51").Value = 7105.67
$ws.Range("C51").Value = '       N/A'
$ws.Range("D51").Value = '       N/A'

$ws.Range("B52").Value = '         N/A'
$ws.Range("C52").Value = '         N/A'
$ws.Range("D52").Value = 158.33

$ws.Range("B53").Value = 7105
$ws.Range("C53").Value = 128
$ws.Range("D53").Value = 158.33

# --- Selection / window state (closest achievable) ---
$ws.Range("D53").Select()

Write-Output "done"
